$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 485) holds the "Förändrad" date serial value.
# Update every cell in that range from 45186 to 45188, keeping formatting intact.
$range = $ws.Range("C2:C485")
$range.Value = 45188
